# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.477.95"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.575.57"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "288.32"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3685"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.80%  "
$ws.Range("E8").Value = "  -3.58%  "
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.151"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07560"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.19%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("E13").Value = "  -0.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.948"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "1.570.76"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001122"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06737"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.391"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "16.53"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.02%  "
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").Value = "22.474.22"
$ws.Range("E24").Value = "  +0.38%  "
$ws.Range("E25").Value = "  +0.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.644"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "150.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.45%  "
$ws.Range("E28").Value = "  +0.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.991"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("D31").Value = "1.748.09"
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.093"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.115"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.994"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.854"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08373"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02468"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2240"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06398"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.295"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.364"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.48"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6277"
$ws.Range("D43").Style = "Normal"
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6115"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.781"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.060"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07227"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.02%  "
